# Apply "hybrid bold + color" quantitative-metric highlighting to specific
# bullet points in the resume, matching the target diff exactly.
#
# For each targeted paragraph, the relevant numeric/percentage/dollar
# substrings are given their own run with Bold + color 2C3E50, while the
# surrounding text is left as plain run(s) - exactly what Word does when
# you select a sub-string of a paragraph and apply character formatting to
# it (it splits runs at the selection boundaries).

$d = $word.ActiveDocument

# BGR-packed integer Word expects for Font.Color matching hex RRGGBB
# 2C3E50 -> R=0x2C G=0x3E B=0x50 -> R + G*256 + B*65536
$HighlightColor = 5258796

function Find-ParagraphIndex($exactText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd("`r", "`n", "`v", "`f")
        if ($t -eq $exactText) {
            return $i
        }
    }
    return -1
}

function Highlight-Metrics($exactText, $metrics) {
    $idx = Find-ParagraphIndex $exactText
    if ($idx -eq -1) {
        Write-Output "WARNING: paragraph not found for text: $exactText"
        return
    }
    $p = $d.Paragraphs.Item($idx)
    $pEnd = $p.Range.End
    $searchStart = $p.Range.Start
    foreach ($m in $metrics) {
        $r = $d.Range($searchStart, $pEnd)
        $found = $r.Find.Execute($m, $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
        if ($found) {
            $r.Font.Bold = 1
            $r.Font.Color = $HighlightColor
            $searchStart = $r.End
        } else {
            Write-Output "WARNING: metric '$m' not found in paragraph $idx"
        }
    }
}

# 1) "• Discovered systematic race coding errors ... accuracy from 23% to 64%"
Highlight-Metrics '• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%' @('23%', '64%')

# 2) "• Achieved 87% prediction accuracy ... of 71%, reducing polling error margins from ±4.2% to ±2.1%"
Highlight-Metrics '• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%' @('87%', '71%', '±4.2%', '±2.1%')

# 3) "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
Highlight-Metrics '• Wrote RFP and analyzed bids from 1,200 vendors for research platform development' @('1,200')

# 4) "...became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+"
Highlight-Metrics '• Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+' @('$400M', '$1B')

# 5) "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
Highlight-Metrics '• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M' @('73.5%', '$4.7M')

# 6) "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" (Key Achievements)
Highlight-Metrics '• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%' @('87%', '71%')

Write-Output "Done"
